$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1) "总计" (Total) summary sheet: shift the existing quarter rows down
#    by one and insert the new "2022-Q4" row at the top of the data
#    (row 2), then append the freed-up last row (2021-Q1) at row 9.
# -------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Walk bottom-up so each row picks up the value that used to be one row
# above it before it gets overwritten.
for ($r = 8; $r -ge 2; $r--) {
    $nr = $r + 1
    $bVal = $totalSheet.Cells.Item($r, 2).Value2
    $cVal = $totalSheet.Cells.Item($r, 3).Value2
    $dVal = $totalSheet.Cells.Item($r, 4).Value2
    $totalSheet.Cells.Item($nr, 2).Value = $bVal
    $totalSheet.Cells.Item($nr, 3).Value = $cVal
    $totalSheet.Cells.Item($nr, 4).Value = $dVal
}

# Row 9 is brand-new - clone the formatting of row 8's A cell (style "2":
# bold / centred / bordered) before writing its final value.
$totalSheet.Cells.Item(8, 1).Copy($totalSheet.Cells.Item(9, 1))
$totalSheet.Cells.Item(9, 1).Value = 7

# New top row: 2022-Q4 figures.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 20
$totalSheet.Cells.Item(2, 4).Value = 0.8100000000000001


# -------------------------------------------------------------------
# 2) New "2022-Q4" detail sheet. Duplicate the "2022-Q3" sheet (so the
#    header row / column-A styling ("s=2": bold, centred, bordered) is
#    carried over identically) right after "总计", rename it, then
#    overwrite its data with the 2022-Q4 fund holdings.
# -------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The copied sheet only has formatted rows 1-4 (header + 3 data rows);
# extend the "s=2" column-A formatting down to row 21 by cloning the
# last formatted data cell (A4) before the values are written.
for ($r = 5; $r -le 21; $r++) {
    $newSheet.Cells.Item(4, 1).Copy($newSheet.Cells.Item($r, 1))
}

$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'006682"
$newSheet.Cells.Item(2,3).Value = "景顺长城中证500指数增强A"
$newSheet.Cells.Item(2,4).Value = "'17.02"
$newSheet.Cells.Item(2,5).Value = "'93.72"
$newSheet.Cells.Item(2,6).Value = "'2.36"
$newSheet.Cells.Item(2,7).Value = "'0.4017"
$newSheet.Cells.Item(2,8).Value = 1
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'000978"
$newSheet.Cells.Item(3,3).Value = "景顺长城量化精选股票"
$newSheet.Cells.Item(3,4).Value = "'7.44"
$newSheet.Cells.Item(3,5).Value = "'94.11"
$newSheet.Cells.Item(3,6).Value = "'2.24"
$newSheet.Cells.Item(3,7).Value = "'0.1667"
$newSheet.Cells.Item(3,8).Value = 2
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'013639"
$newSheet.Cells.Item(4,3).Value = "光大保德信中证500指数增强A"
$newSheet.Cells.Item(4,4).Value = "'3.06"
$newSheet.Cells.Item(4,5).Value = "'90.95"
$newSheet.Cells.Item(4,6).Value = "'1.68"
$newSheet.Cells.Item(4,7).Value = "'0.0514"
$newSheet.Cells.Item(4,8).Value = 3
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "'009992"
$newSheet.Cells.Item(5,3).Value = "景顺长城量化成长演化混合"
$newSheet.Cells.Item(5,4).Value = "'2.31"
$newSheet.Cells.Item(5,5).Value = "'91.65"
$newSheet.Cells.Item(5,6).Value = "'2.12"
$newSheet.Cells.Item(5,7).Value = "'0.0490"
$newSheet.Cells.Item(5,8).Value = 10
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "'162216"
$newSheet.Cells.Item(6,3).Value = "泰达宏利中证500指数增强（LOF）"
$newSheet.Cells.Item(6,4).Value = "'2.78"
$newSheet.Cells.Item(6,5).Value = "'93.77"
$newSheet.Cells.Item(6,6).Value = "'1.38"
$newSheet.Cells.Item(6,7).Value = "'0.0384"
$newSheet.Cells.Item(6,8).Value = 1
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "'012080"
$newSheet.Cells.Item(7,3).Value = "易方达中证500指数量化增强A"
$newSheet.Cells.Item(7,4).Value = "'3.80"
$newSheet.Cells.Item(7,5).Value = "'93.92"
$newSheet.Cells.Item(7,6).Value = "'0.98"
$newSheet.Cells.Item(7,7).Value = "'0.0372"
$newSheet.Cells.Item(7,8).Value = 7
$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "'510200"
$newSheet.Cells.Item(8,3).Value = "汇安上证证券ETF"
$newSheet.Cells.Item(8,4).Value = "'0.65"
$newSheet.Cells.Item(8,5).Value = "'95.50"
$newSheet.Cells.Item(8,6).Value = "'3.07"
$newSheet.Cells.Item(8,7).Value = "'0.0200"
$newSheet.Cells.Item(8,8).Value = 9
$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "'012081"
$newSheet.Cells.Item(9,3).Value = "易方达中证500指数量化增强C"
$newSheet.Cells.Item(9,4).Value = "'1.30"
$newSheet.Cells.Item(9,5).Value = "'93.92"
$newSheet.Cells.Item(9,6).Value = "'0.98"
$newSheet.Cells.Item(9,7).Value = "'0.0127"
$newSheet.Cells.Item(9,8).Value = 7
$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "'013640"
$newSheet.Cells.Item(10,3).Value = "光大保德信中证500指数增强C"
$newSheet.Cells.Item(10,4).Value = "'0.70"
$newSheet.Cells.Item(10,5).Value = "'90.95"
$newSheet.Cells.Item(10,6).Value = "'1.68"
$newSheet.Cells.Item(10,7).Value = "'0.0118"
$newSheet.Cells.Item(10,8).Value = 3
$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = "'003760"
$newSheet.Cells.Item(11,3).Value = "国泰中证500指数增强A"
$newSheet.Cells.Item(11,4).Value = "'0.48"
$newSheet.Cells.Item(11,5).Value = "'91.42"
$newSheet.Cells.Item(11,6).Value = "'1.78"
$newSheet.Cells.Item(11,7).Value = "'0.0085"
$newSheet.Cells.Item(11,8).Value = 5
$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).Value = "'005260"
$newSheet.Cells.Item(12,3).Value = "银华稳健增利灵活配置混合A"
$newSheet.Cells.Item(12,4).Value = "'0.33"
$newSheet.Cells.Item(12,5).Value = "'91.18"
$newSheet.Cells.Item(12,6).Value = "'0.79"
$newSheet.Cells.Item(12,7).Value = "'0.0026"
$newSheet.Cells.Item(12,8).Value = 2
$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).Value = "'012926"
$newSheet.Cells.Item(13,3).Value = "民生加银中证500指数增强A"
$newSheet.Cells.Item(13,4).Value = "'0.20"
$newSheet.Cells.Item(13,5).Value = "'87.82"
$newSheet.Cells.Item(13,6).Value = "'0.86"
$newSheet.Cells.Item(13,7).Value = "'0.0017"
$newSheet.Cells.Item(13,8).Value = 7
$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).Value = "'005261"
$newSheet.Cells.Item(14,3).Value = "银华稳健增利灵活配置混合C"
$newSheet.Cells.Item(14,4).Value = "'0.21"
$newSheet.Cells.Item(14,5).Value = "'91.18"
$newSheet.Cells.Item(14,6).Value = "'0.79"
$newSheet.Cells.Item(14,7).Value = "'0.0017"
$newSheet.Cells.Item(14,8).Value = 2
$newSheet.Cells.Item(15,1).Value = 13
$newSheet.Cells.Item(15,2).Value = "'005966"
$newSheet.Cells.Item(15,3).Value = "安信中证500指数增强C"
$newSheet.Cells.Item(15,4).Value = "'0.16"
$newSheet.Cells.Item(15,5).Value = "'88.79"
$newSheet.Cells.Item(15,6).Value = "'0.87"
$newSheet.Cells.Item(15,7).Value = "'0.0014"
$newSheet.Cells.Item(15,8).Value = 10
$newSheet.Cells.Item(16,1).Value = 14
$newSheet.Cells.Item(16,2).Value = "'012927"
$newSheet.Cells.Item(16,3).Value = "民生加银中证500指数增强C"
$newSheet.Cells.Item(16,4).Value = "'0.14"
$newSheet.Cells.Item(16,5).Value = "'87.82"
$newSheet.Cells.Item(16,6).Value = "'0.86"
$newSheet.Cells.Item(16,7).Value = "'0.0012"
$newSheet.Cells.Item(16,8).Value = 7
$newSheet.Cells.Item(17,1).Value = 15
$newSheet.Cells.Item(17,2).Value = "'006783"
$newSheet.Cells.Item(17,3).Value = "红土创新中证500指数增强A"
$newSheet.Cells.Item(17,4).Value = "'0.05"
$newSheet.Cells.Item(17,5).Value = "'92.80"
$newSheet.Cells.Item(17,6).Value = "'2.23"
$newSheet.Cells.Item(17,7).Value = "'0.0011"
$newSheet.Cells.Item(17,8).Value = 3
$newSheet.Cells.Item(18,1).Value = 16
$newSheet.Cells.Item(18,2).Value = "'005965"
$newSheet.Cells.Item(18,3).Value = "安信中证500指数增强A"
$newSheet.Cells.Item(18,4).Value = "'0.12"
$newSheet.Cells.Item(18,5).Value = "'88.79"
$newSheet.Cells.Item(18,6).Value = "'0.87"
$newSheet.Cells.Item(18,7).Value = "'0.0010"
$newSheet.Cells.Item(18,8).Value = 10
$newSheet.Cells.Item(19,1).Value = 17
$newSheet.Cells.Item(19,2).Value = "'006784"
$newSheet.Cells.Item(19,3).Value = "红土创新中证500指数增强C"
$newSheet.Cells.Item(19,4).Value = "'0.04"
$newSheet.Cells.Item(19,5).Value = "'92.80"
$newSheet.Cells.Item(19,6).Value = "'2.23"
$newSheet.Cells.Item(19,7).Value = "'0.0009"
$newSheet.Cells.Item(19,8).Value = 3
$newSheet.Cells.Item(20,1).Value = 18
$newSheet.Cells.Item(20,2).Value = "'003761"
$newSheet.Cells.Item(20,3).Value = "国泰中证500指数增强C"
$newSheet.Cells.Item(20,4).Value = "'0.04"
$newSheet.Cells.Item(20,5).Value = "'91.42"
$newSheet.Cells.Item(20,6).Value = "'1.78"
$newSheet.Cells.Item(20,7).Value = "'0.0007"
$newSheet.Cells.Item(20,8).Value = 5
$newSheet.Cells.Item(21,1).Value = 19
$newSheet.Cells.Item(21,2).Value = "'016935"
$newSheet.Cells.Item(21,3).Value = "景顺长城中证500指数增强C"
$newSheet.Cells.Item(21,4).Value = "'0.00"
$newSheet.Cells.Item(21,5).Value = "'93.72"
$newSheet.Cells.Item(21,6).Value = "'2.36"
$newSheet.Cells.Item(21,7).Value = 0
$newSheet.Cells.Item(21,8).Value = 1

# Restore the originally active sheet/tab ("总计") - copying/renaming
# sheets above shifts Excel's active-sheet cursor to the newly created
# sheet, which the source workbook did not have selected.
$totalSheet.Activate()
